$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The database of workers/periods was updated: previous account statement (EC)
# rows were removed and new ones added.
#   Row 16: now BLADIMIR PACHECO AYALA (20351521), period 2203
#   Row 17: now CARLOS HUMBERTO MALDONADO RANGEL (1092345081), period 2204
#   Row 18: now CARLOS HUMBERTO MALDONADO RANGEL (1092345081), period 2203
#   Row 19: unchanged - LUIS DAVID ORTIZ CALA (1098778032), period 2309

$ws.Range("C16").Value = "20351521"
$ws.Range("D16").Value = "BLADIMIR PACHECO AYALA"
$ws.Range("E16").Value = "2203"

$ws.Range("C17").Value = "1092345081"
$ws.Range("D17").Value = "CARLOS HUMBERTO MALDONADO RANGEL"
$ws.Range("E17").Value = "2204"

$ws.Range("C18").Value = "1092345081"
$ws.Range("D18").Value = "CARLOS HUMBERTO MALDONADO RANGEL"
$ws.Range("E18").Value = "2203"
